$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1").Style = $ws.Range("H1").Style
$ws.Range("J1").Style = $ws.Range("H1").Style

# Data values for columns I (I0) and J (IF), rows 2-74
$iVals = @(9,7,3,8,5,7,5,9,7,7,6,6,6,6,7,4,6,10,7,8,8,8,9,9,9,8,10,7,8,7,8,6,8,8,8,7,7,7,7,8,8,8,7,8,8,8,7,8,8,7,7,8,8,9,8,8,8,7,7,6,8,6,9,6,5,6,8,1,6,7,7,6,7)
$jVals = @(9,7,4,8,6,8,5,9,7,8,7,7,6,7,8,5,6,10,7,8,8,8,9,10,9,8,10,8,8,7,8,7,8,8,8,8,7,7,8,8,8,8,7,8,8,8,7,8,8,8,8,8,8,9,8,8,8,8,7,7,8,6,9,6,5,7,8,1,6,7,7,6,7)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
